$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "[-, 'MCT-3A-Manut. Mecânica', -, -]"
$ws.Range("E2").Value = "[-, 'MCT-3A-Manut. Mecânica', -, -]"
$ws.Range("D3").Value = "[-, 'MCT-3A-Manut. Mecânica', -, -]"
$ws.Range("F3").Value = "-"
$ws.Range("F4").Value = "-"
$ws.Range("F6").Value = "-"
$ws.Range("F7").Value = "-"
$ws.Range("E8").Value = "[-, -, 'MCT-3A-Manut. Mecânica', -]"
